$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the blank "spacer" paragraphs that sit between every heading and
#    body paragraph. We keep three of them (the ones that remain blank in the
#    target document): the one right after "## De la Gran Reunificacion", and
#    the two near the signature block at the very end.
#    Deleting a paragraph's Range merges it away (paragraph mark included),
#    which is exactly what the diff shows happening everywhere else.
#    Indices are processed from highest to lowest so earlier ones stay valid.
# ---------------------------------------------------------------------------
$blankParasToRemove = @(35, 33, 29, 26, 24, 22, 20, 18, 16, 14, 12, 10, 8, 6, 4, 2)
foreach ($i in $blankParasToRemove) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) Strip the leading Markdown-style "# " / "## " markers from the title and
#    section headings now that they are plain Word headings.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("# Crónicas del Imperio ", $true, $false, $false, $false, $false, $true, 1, $false, "Crónicas del Imperio ", 2) | Out-Null
$d.Content.Find.Execute("## De los Orígenes del Imperio y sus Padres Fundadores", $true, $false, $false, $false, $false, $true, 1, $false, "De los Orígenes del Imperio y sus Padres Fundadores", 2) | Out-Null
$d.Content.Find.Execute("## Del Himno Imperial", $true, $false, $false, $false, $false, $true, 1, $false, "Del Himno Imperial", 2) | Out-Null
$d.Content.Find.Execute("## De la Economía y Moneda Imperial", $true, $false, $false, $false, $false, $true, 1, $false, "De la Economía y Moneda Imperial", 2) | Out-Null
$d.Content.Find.Execute("## De las Costumbres y Rituales del Imperio", $true, $false, $false, $false, $false, $true, 1, $false, "De las Costumbres y Rituales del Imperio", 2) | Out-Null
$d.Content.Find.Execute("## De la Gran Reunificación", $true, $false, $false, $false, $false, $true, 1, $false, "De la Gran Reunificación", 2) | Out-Null
$d.Content.Find.Execute("## De las Profecías y el Futuro", $true, $false, $false, $false, $false, $true, 1, $false, "De las Profecías y el Futuro", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Justify (both-align) every paragraph in the document.
# ---------------------------------------------------------------------------
$d.Paragraphs.Alignment = 3

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
